# "finish the last benchmarks"
# Append three new HC6(...) benchmark rows (18-20) to the "5" worksheet's
# M:W results table, and make that sheet the active / selected one
# (M25 selection) instead of the "27" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5")

# --- new benchmark rows ------------------------------------------------
# row 18: HC6(11)
$ws.Cells.Item(18, 13).Value = "HC6(11)"   # M18
$ws.Cells.Item(18, 14).Value = 32          # N18
$ws.Cells.Item(18, 15).Value = 18          # O18
$ws.Cells.Item(18, 16).Formula = "=N18-O18"            # P18
$ws.Cells.Item(18, 17).Value = 14          # Q18
$ws.Cells.Item(18, 18).Value = 20          # R18
$ws.Cells.Item(18, 19).Formula = "=(O18-Q18)/N18"       # S18
$ws.Cells.Item(18, 20).Formula = "=(P18-R18+Q18)/N18"   # T18
$ws.Cells.Item(18, 21).Formula = "= 1 -R18/N18"          # U18
$ws.Cells.Item(18, 22).Value = 1           # V18
$ws.Cells.Item(18, 23).Formula = "=V18/N18"              # W18

# row 19: HC6(13)
$ws.Cells.Item(19, 13).Value = "HC6(13)"   # M19
$ws.Cells.Item(19, 14).Value = 32          # N19
$ws.Cells.Item(19, 15).Value = 18          # O19
$ws.Cells.Item(19, 16).Formula = "=N19-O19"            # P19
$ws.Cells.Item(19, 17).Value = 15          # Q19
$ws.Cells.Item(19, 18).Value = 21          # R19
$ws.Cells.Item(19, 19).Formula = "=(O19-Q19)/N19"       # S19
$ws.Cells.Item(19, 20).Formula = "=(P19-R19+Q19)/N19"   # T19
$ws.Cells.Item(19, 21).Formula = "= 1 -R19/N19"          # U19
$ws.Cells.Item(19, 22).Value = 0           # V19
$ws.Cells.Item(19, 23).Formula = "=V19/N19"              # W19

# row 20: HC6(15)
$ws.Cells.Item(20, 13).Value = "HC6(15)"   # M20
$ws.Cells.Item(20, 14).Value = 32          # N20
$ws.Cells.Item(20, 15).Value = 18          # O20
$ws.Cells.Item(20, 16).Formula = "=N20-O20"            # P20
$ws.Cells.Item(20, 17).Value = 18          # Q20
$ws.Cells.Item(20, 18).Value = 26          # R20
$ws.Cells.Item(20, 19).Formula = "=(O20-Q20)/N20"       # S20
$ws.Cells.Item(20, 20).Formula = "=(P20-R20+Q20)/N20"   # T20
$ws.Cells.Item(20, 21).Formula = "= 1 -R20/N20"          # U20
$ws.Cells.Item(20, 22).Value = 1           # V20
$ws.Cells.Item(20, 23).Formula = "=V20/N20"              # W20

# --- selection / active sheet ------------------------------------------
# Move the active tab from "27" to "5", and leave the cursor on M25.
$ws.Activate()
$ws.Range("M25").Select()
